# Auto-generated script applying scheduled market-data refresh to all 8 Leve sheets.
# Each (row, column) below corresponds to currentAveragePrice* / LeveProfit* figures
# pulled by the scheduled runner; values are written verbatim, and cells that the
# refresh leaves blank are cleared so the workbook matches the upstream snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2095.8
$ws.Range("I12").Value = 2095.8
$ws.Range("K12").Value = 2095.8
$ws.Range("M12").Value = -1925.8
$ws.Range("H29").Value = 998.6667
$ws.Range("I29").Value = 999
$ws.Range("J29").Value = 998
$ws.Range("K29").Value = 2997
$ws.Range("L29").Value = 2994
$ws.Range("M29").Value = -2716
$ws.Range("N29").Value = -3556
$ws.Range("H43").Value = 1225
$ws.Range("I43").Value = 1100
$ws.Range("J43").Value = 1350
$ws.Range("K43").Value = 1100
$ws.Range("L43").Value = 1350
$ws.Range("M43").Value = -1031
$ws.Range("N43").Value = -1488
$ws.Range("H55").Value = 1017.94116
$ws.Range("I55").Value = 1595.25
$ws.Range("J55").Value = 504.77777
$ws.Range("K55").Value = 1595.25
$ws.Range("L55").Value = 504.77777
$ws.Range("M55").Value = -1381.25
$ws.Range("N55").Value = -932.7777699999999
$ws.Range("H76").Value = 50006000
$ws.Range("J76").Value = 12000
$ws.Range("L76").Value = 12000
$ws.Range("N76").Value = -12630
$ws.Range("H79").Value = 50006000
$ws.Range("J79").Value = 12000
$ws.Range("L79").Value = 12000
$ws.Range("N79").Value = -14184
$ws.Range("H125").Value = 633.8570999999999
$ws.Range("J125").Value = 568.4
$ws.Range("L125").Value = 5115.599999999999
$ws.Range("N125").Value = -10035.6
$ws.Range("H127").Value = 1878
$ws.Range("I127").Value = 1878
$ws.Range("K127").Value = 5634
$ws.Range("M127").Value = -674
$ws.Range("H131").Value = 1166.7142
$ws.Range("I131").Value = 1166.7142
$ws.Range("K131").Value = 3500.1426
$ws.Range("M131").Value = 1539.8574
$ws.Range("H138").Value = 3048
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4511.6
$ws.Range("I32").Value = 4735.1113
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 4735.1113
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = -4448.1113
$ws.Range("N32").Value = -3074
$ws.Range("H74").Value = 1472
$ws.Range("I74").Value = 1472
$ws.Range("K74").Value = 1472
$ws.Range("M74").Value = -598
$ws.Range("H76").Value = 47495
$ws.Range("J76").Value = 47495
$ws.Range("L76").Value = 47495
$ws.Range("N76").Value = -48171
$ws.Range("H77").Value = 1472
$ws.Range("I77").Value = 1472
$ws.Range("K77").Value = 7360
$ws.Range("M77").Value = -2992
$ws.Range("H79").Value = 47495
$ws.Range("J79").Value = 47495
$ws.Range("L79").Value = 47495
$ws.Range("N79").Value = -49835
$ws.Range("H97").Value = 1526.7778
$ws.Range("I97").Value = 1028.3529
$ws.Range("K97").Value = 1028.3529
$ws.Range("M97").Value = -532.3529000000001
$ws.Range("H130").Value = 21096.4
$ws.Range("I130").Value = 20815
$ws.Range("K130").Value = 20815
$ws.Range("M130").Value = -15795
$ws.Range("H132").Value = 6997.222
$ws.Range("I132").Value = 7496
$ws.Range("J132").Value = 5999.6665
$ws.Range("K132").Value = 22488
$ws.Range("L132").Value = 17998.9995
$ws.Range("M132").Value = -19958
$ws.Range("N132").Value = -23058.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2215
$ws.Range("I20").Value = 2418.25
$ws.Range("J20").Value = 2011.75
$ws.Range("K20").Value = 2418.25
$ws.Range("L20").Value = 2011.75
$ws.Range("M20").Value = -2171.25
$ws.Range("N20").Value = -2505.75
$ws.Range("H86").Value = 8555.24
$ws.Range("I86").Value = 7946.579
$ws.Range("J86").Value = 10482.667
$ws.Range("K86").Value = 7946.579
$ws.Range("L86").Value = 10482.667
$ws.Range("M86").Value = -6823.579
$ws.Range("N86").Value = -12728.667
$ws.Range("H89").Value = 8555.24
$ws.Range("I89").Value = 7946.579
$ws.Range("J89").Value = 10482.667
$ws.Range("K89").Value = 39732.895
$ws.Range("L89").Value = 52413.335
$ws.Range("M89").Value = -34116.895
$ws.Range("N89").Value = -63645.335
$ws.Range("H94").Value = 718.35297
$ws.Range("I94").Value = 753.875
$ws.Range("K94").Value = 753.875
$ws.Range("M94").Value = -302.875
$ws.Range("H105").Value = 4124.8335
$ws.Range("I105").Value = 4349.8
$ws.Range("K105").Value = 4349.8
$ws.Range("M105").Value = -2602.8
$ws.Range("H134").Value = 8275.299999999999
$ws.Range("I134").Value = 8230.375
$ws.Range("K134").Value = 24691.125
$ws.Range("M134").Value = -22156.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H25").Value = 24999.5
$ws.Range("I25").Value = 24999
$ws.Range("K25").Value = 24999
$ws.Range("M25").Value = -24825
$ws.Range("H31").Value = 1873.05
$ws.Range("I31").Value = 1711
$ws.Range("K31").Value = 1711
$ws.Range("M31").Value = -1416
$ws.Range("H34").Value = 1873.05
$ws.Range("I34").Value = 1711
$ws.Range("K34").Value = 1711
$ws.Range("M34").Value = -1509
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H105").Value = 3230.5
$ws.Range("I105").Value = 2450
$ws.Range("K105").Value = 2450
$ws.Range("M105").Value = -703

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 1999.3334
$ws.Range("J123").Value = 498
$ws.Range("L123").Value = 1494
$ws.Range("N123").Value = -6394
$ws.Range("H137").Value = 4750
$ws.Range("I137").Value = 4750
$ws.Range("K137").Value = 14250
$ws.Range("M137").Value = -9150
$ws.Range("H141").Value = 7776.8887
$ws.Range("I141").Value = 7776.8887
$ws.Range("K141").Value = 23330.6661
$ws.Range("M141").Value = -18150.6661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 25800
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 25800
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 25800
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -26184
$ws.Range("H80").Value = 1941.4166
$ws.Range("I80").Value = 1729.5714
$ws.Range("J80").Value = 2238
$ws.Range("K80").Value = 1729.5714
$ws.Range("L80").Value = 2238
$ws.Range("M80").Value = -731.5714
$ws.Range("N80").Value = -4234
$ws.Range("H83").Value = 1941.4166
$ws.Range("I83").Value = 1729.5714
$ws.Range("J83").Value = 2238
$ws.Range("K83").Value = 8647.857
$ws.Range("L83").Value = 11190
$ws.Range("M83").Value = -3655.857
$ws.Range("N83").Value = -21174
$ws.Range("H97").Value = 923.375
$ws.Range("I97").Value = 923.375
$ws.Range("K97").Value = 923.375
$ws.Range("M97").Value = -427.375
$ws.Range("H102").Value = 2663.0527
$ws.Range("I102").Value = 2543.6875
$ws.Range("J102").Value = 3299.6667
$ws.Range("K102").Value = 2543.6875
$ws.Range("L102").Value = 3299.6667
$ws.Range("M102").Value = -921.6875
$ws.Range("N102").Value = -6543.6667
$ws.Range("H122").Value = 5025.3335
$ws.Range("I122").Value = 4388
$ws.Range("K122").Value = 13164
$ws.Range("M122").Value = -10714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5077.143
$ws.Range("I7").Value = 5168.143
$ws.Range("J7").Value = 4986.143
$ws.Range("K7").Value = 5168.143
$ws.Range("L7").Value = 4986.143
$ws.Range("M7").Value = -5056.143
$ws.Range("N7").Value = -5210.143
$ws.Range("H40").Value = 2322.1667
$ws.Range("I40").Value = 1983.5
$ws.Range("K40").Value = 1983.5
$ws.Range("M40").Value = -1847.5
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("N44").Value = -15912
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H122").Value = 3316.3635
$ws.Range("I122").Value = 3164.5557
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 9493.667099999999
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -7043.667099999999
$ws.Range("N122").Value = -16898.5
$ws.Range("H126").Value = 5077.143
$ws.Range("I126").Value = 5168.143
$ws.Range("J126").Value = 4986.143
$ws.Range("K126").Value = 15504.429
$ws.Range("L126").Value = 14958.429
$ws.Range("M126").Value = -13034.429
$ws.Range("N126").Value = -19898.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5648.25
$ws.Range("I122").Value = 4888.5454
$ws.Range("J122").Value = 14005
$ws.Range("K122").Value = 14665.6362
$ws.Range("L122").Value = 42015
$ws.Range("M122").Value = -12215.6362
$ws.Range("N122").Value = -46915
$ws.Range("H126").Value = 1919.1
$ws.Range("I126").Value = 1354.5555
$ws.Range("K126").Value = 4063.6665
$ws.Range("M126").Value = -1593.6665
